$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$old1 = "✅ 1000 Bs = 4.53 = 17994.97 pesos"
$new1 = "✅ 1000 Bs = 4.62 = 18418.01 pesos"
$old2 = "✅ 17994.97 pesos = 4.51 = 951.15 Bs"
$new2 = "✅ 18418.01 pesos = 4.61 = 964.27 Bs"

$text = [string]$wsHoja1.Range("A1").Value2
$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$wsHoja1.Range("A1").Value2 = $text

# --- Update the tasas rates on tasas!N10, O10, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 216.5
$wsTasas.Range("O10").Value = 3987.5
$wsTasas.Range("O12").Value = 209
